$wb = $excel.ActiveWorkbook

# --- "PCiCDTdtTDM" sheet ---------------------------------------------------
# We now allow for twice the potential identified in the BLUE Shifts
# scenario, so every lever formula is doubled.
$wsPCi = $wb.Worksheets.Item("PCiCDTdtTDM")
$wsPCi.Range("B2").Formula = "=Calcs!B5*2"
$wsPCi.Range("B3").Formula = "=Calcs!C5*2"
$wsPCi.Range("C3").Formula = "=Calcs!B11*2"
$wsPCi.Range("B4").Formula = "=Calcs!D5*2"
$wsPCi.Range("B5").Formula = "=Calcs!E5*2"
$wsPCi.Range("C5").Formula = "=Calcs!C11*2"
$wsPCi.Range("B6").Formula = "=Calcs!F5*2"
$wsPCi.Range("B7").Formula = "=Calcs!G5*2"

[void]$wsPCi.Range("C6").Select()

# --- "About" sheet -------------------------------------------------------
# Insert a new note row (plus a following blank spacer row, matching the
# existing blank-row rhythm of the sheet) right before the old "While data
# is given for 2050..." paragraph, pushing everything below down by two
# rows.
$wsAbout = $wb.Worksheets.Item("About")
[void]$wsAbout.Rows("18:19").Insert()
$wsAbout.Range("A18").Value = "We allow for twice the potential identified in the BLUE Shifts scenario."

# Leave the selection where the new text was typed (whole-row selection on
# the blank spacer row beneath it), like Excel does after inserting rows.
# Selecting this sheet last also keeps "About" as the active tab, matching
# the original workbook.
[void]$wsAbout.Rows("19:19").Select()
